$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.083.65"
$ws.Range("E2").Value = "  +0.65%  "

$ws.Range("D3").Value = "3.168.48"
$ws.Range("E3").Value = "  +4.03%  "

$ws.Range("E4").Value = "  -0.61%  "

$ws.Range("D5").Value = "'216.31"
$ws.Range("E5").Value = "  +1.80%  "

$ws.Range("D6").Value = "'629.54"
$ws.Range("E6").Value = "  +2.78%  "

$ws.Range("E7").Value = "  +32.77%  "

$ws.Range("E8").Value = "  +0.82%  "

$ws.Range("E9").Value = "  -0.10%  "

$ws.Range("D10").Value = "3.164.95"
$ws.Range("E10").Value = "  +4.06%  "

$ws.Range("D11").Value = "'0.762"
$ws.Range("E11").Value = "  +12.60%  "

$ws.Range("E12").Value = "  +8.23%  "

$ws.Range("D13").Value = "'5.80"
$ws.Range("E13").Value = "  +8.71%  "

$ws.Range("E14").Value = "  +2.24%  "

$ws.Range("D15").Value = "'35.19"
$ws.Range("E15").Value = "  +8.32%  "

$ws.Range("D16").Value = "90.832.17"
$ws.Range("E16").Value = "  +0.59%  "

$ws.Range("D17").Value = "3.753.01"
$ws.Range("E17").Value = "  +3.92%  "

$ws.Range("D18").Value = "3.174.37"
$ws.Range("E18").Value = "  +6.12%  "

$ws.Range("D19").Value = "'3.78"
$ws.Range("E19").Value = "  +12.53%  "

$ws.Range("D20").Value = "'14.74"
$ws.Range("E20").Value = "  +9.33%  "

$ws.Range("D21").Value = "'481.52"
$ws.Range("E21").Value = "  +12.93%  "

$ws.Range("E22").Value = "  -4.07%  "

$ws.Range("E23").Value = "  +10.40%  "

$ws.Range("D24").Value = "'5.17"
$ws.Range("E24").Value = "  +2.46%  "

$ws.Range("D25").Value = "'97.20"
$ws.Range("E25").Value = "  +16.95%  "

$ws.Range("D26").Value = "'5.93"
$ws.Range("E26").Value = "  +10.57%  "

$ws.Range("D27").Value = "'12.37"
$ws.Range("E27").Value = "  +6.48%  "

$ws.Range("D28").Value = "3.340.43"
$ws.Range("E28").Value = "  +4.19%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("D30").Value = "'0.163"
$ws.Range("E30").Value = "  +0.79%  "

$ws.Range("D31").Value = "'9.33"
$ws.Range("E31").Value = "  +8.08%  "

$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  -0.25%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'28.23"
$ws.Range("E33").Value = "  +23.77%  "

$ws.Range("E34").Value = "  +44.07%  "

$ws.Range("D35").Value = "'526.02"
$ws.Range("E35").Value = "  +4.46%  "

$ws.Range("E36").Value = "  +7.22%  "

$ws.Range("E37").Value = "  +9.57%  "

$ws.Range("D38").Value = "'3.64"
$ws.Range("E38").Value = "  -3.54%  "

$ws.Range("D39").Value = "'6.99"
$ws.Range("E39").Value = "  +4.53%  "

$ws.Range("E40").Value = "  +4.68%  "

$ws.Range("D41").Value = "'0.0927"
$ws.Range("E41").Value = "  +32.99%  "

$ws.Range("E42").Value = "  -0.22%  "

$ws.Range("E43").Value = "  +17.31%  "

$ws.Range("E44").Value = "  -0.08%  "

$ws.Range("E45").Value = "  +8.68%  "

$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("D47").Value = "'0.717"
$ws.Range("E47").Value = "  +20.86%  "

$ws.Range("D48").Value = "'4.67"
$ws.Range("E48").Value = "  +11.11%  "

$ws.Range("D49").Value = "'150.41"
$ws.Range("E49").Value = "  +4.84%  "

$ws.Range("E50").Value = "  +12.58%  "

$ws.Range("D51").Value = "'45.58"
$ws.Range("E51").Value = "  +4.72%  "
